# Lume.xlsx: "some bug fixes" -- add an Hour/Hyperminute reference table and
# a timer-frequency calculator block to Sheet1, drop the stray scratch
# formula in G15, and make Sheet1 (instead of Money) the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- remove the stray one-off calculation that used to live in G15 -------
$ws.Range("G15").ClearContents()

# --- Hour / Hyperminute lookup table (rows 20-32, cols D/E and G/H/I) ----
$ws.Range("D20").Value = "Hour"
$ws.Range("G20").Value = "Hyperminute"

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "12 и 11"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "1 и 12"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "2 и 1"
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = "3 и 2"
$ws.Range("D25").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("D27").Value = 6
$ws.Range("D28").Value = 7
$ws.Range("D29").Value = 8
$ws.Range("D30").Value = 9
$ws.Range("D31").Value = 10
$ws.Range("E31").Value = "10 и 9"
$ws.Range("D32").Value = 11
$ws.Range("E32").Value = "11 и 10"

$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "12 и 11.5"
$ws.Range("G22").Formula = "=G21+1"
$ws.Range("H22").Value = "0.5 и 0"
$ws.Range("I22").Value = "!"
$ws.Range("G23").Formula = "=G22+1"
$ws.Range("H23").Value = "1 и 0.5"
$ws.Range("G24:G44").Formula = "=G23+1"

# --- timer-frequency calculator (rows 32-41, cols A-C) --------------------
$ws.Range("A32").Value = "F CPU"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "MHz"

$ws.Range("A33").Value = "Timer ovf freq"
$ws.Range("B33").Formula = "=B32*1000000/256"
$ws.Range("C33").Value = "Hz"

$ws.Range("A38").Value = "Divisor"
$ws.Range("B38").Value = 1

$ws.Range("A39").Value = "Timer input freq"
$ws.Range("B39").Formula = "=B32*1000000/B38"
$ws.Range("C39").Value = "Hz"

$ws.Range("A40").Value = "ICR"
$ws.Range("B40").Value = 255

$ws.Range("A41").Value = "OVF freq"
$ws.Range("B41").Formula = "=B39/B40"
$ws.Range("C41").Value = "Hz"

# --- widen column E a bit so the hour labels aren't clipped ---------------
$ws.Columns.Item(5).ColumnWidth = 9.5

# --- make Sheet1 the active sheet/tab (was "Money") ------------------------
$ws.Activate()
$ws.Range("C44").Select()
